# Insert a new weekly price-report row for Cilantro (Vega Modelo de Temuco)
# at row 217, shifting all subsequent rows (217-253) down by one (218-254).
# This mirrors a new observation being added at the top of the recent-dates
# block while the rest of the historical rows keep their relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 217 (and everything below it) down by one row.
$ws.Rows.Item(217).Insert()

# Populate the newly inserted row 217 with the new observation.
$ws.Cells.Item(217, 1).Value  = 10
$ws.Cells.Item(217, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(217, 3).Value  = "La Araucanía"
$ws.Cells.Item(217, 4).Value  = 44504
$ws.Cells.Item(217, 5).Value  = 9
$ws.Cells.Item(217, 6).Value  = 100112040
$ws.Cells.Item(217, 7).Value  = "Cilantro"
$ws.Cells.Item(217, 8).Value  = "Sin especificar"
$ws.Cells.Item(217, 9).Value  = "Primera"
$ws.Cells.Item(217, 10).Value = 110
$ws.Cells.Item(217, 11).Value = 4500
$ws.Cells.Item(217, 12).Value = 4500
$ws.Cells.Item(217, 13).Value = 4500
$ws.Cells.Item(217, 14).Value = "`$/docena de atados (2 kilos)"
$ws.Cells.Item(217, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(217, 16).Value = 2250
$ws.Cells.Item(217, 17).Value = 2
$ws.Cells.Item(217, 18).Value = "Hortaliza"
